# Automatische test-sync: 2025-08-03 18:17:50
#
# Adds a new "Testmail #3" row (row 31) to the Logs sheet, extends the
# conditional-formatting ranges that cover the data rows to include the
# new row, and updates the Dashboard pivot/summary sheet so the
# "Planning / Afspraak" category count (now 9) sorts above "Overig".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Logs sheet: append the new log entry in row 31
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A31").Value = "Kun jij dit afhandelen?"
$logs.Range("B31").Value = "mailmind.test@zohomail.eu"
$logs.Range("C31").Value = "Testmail #3: Kun jij dit afhandelen?"
$logs.Range("D31").Value = "Planning / Afspraak"
$logs.Range("E31").Value = "Bedankt, we hebben dit doorgestuurd naar planning@bedrijf.nl."
$logs.Range("F31").Value = "2025-08-03 18:17:02"
$logs.Range("G31").Value = "Ja"
$logs.Range("H31").Value = "Ja"
$logs.Range("I31").Value = "Nee"
$logs.Range("J31").Value = "Nee"

# ---------------------------------------------------------------------
# 2) Grow the conditional-formatting ranges from row 30 to row 31 so the
#    new row keeps getting highlighted like the rest of the table.
# ---------------------------------------------------------------------
$cfColumns = @("D", "G", "H", "I", "J")
foreach ($col in $cfColumns) {
    $oldRange = $logs.Range("$($col)2:$($col)30")
    $newRange = $logs.Range("$($col)2:$($col)31")
    $conditions = $oldRange.FormatConditions
    for ($i = 1; $i -le $conditions.Count; $i++) {
        $conditions.Item($i).ModifyAppliesToRange($newRange)
    }
}

# ---------------------------------------------------------------------
# 3) Dashboard sheet: "Planning / Afspraak" now ties "Overig" at 9, and
#    moves above it in the summary table.
# ---------------------------------------------------------------------
$dashboard = $wb.Worksheets.Item("Dashboard")

$dashboard.Range("A2").Value = "Planning / Afspraak"
$dashboard.Range("B2").Value = 9

$dashboard.Range("A3").Value = "Overig"
$dashboard.Range("B3").Value = 9
